# spring 24 reg season complete
# Update matchup average values on the "Nine" worksheet to reflect the
# latest regular-season results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 10.31
$ws.Range("E3").Value = 10.76

$ws.Range("C4").Value = 9.69
$ws.Range("E4").Value = 10.64
$ws.Range("G4").Value = 10.12

$ws.Range("C5").Value = 9.19
$ws.Range("D5").Value = 9.359999999999999
$ws.Range("F5").Value = 10.22
$ws.Range("G5").Value = 9.68

$ws.Range("E6").Value = 9.779999999999999
$ws.Range("G6").Value = 10.31
$ws.Range("H6").Value = 10.4

$ws.Range("D7").Value = 10.06
$ws.Range("E7").Value = 10.32
$ws.Range("F7").Value = 9.69
$ws.Range("I7").Value = 6.58

$ws.Range("F8").Value = 9.6
$ws.Range("J8").Value = 11.22

$ws.Range("G9").Value = 13.42

$ws.Range("H10").Value = 8.779999999999999
